$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 47.21827566666667
$ws.Range("H2").Value = 141.654827
$ws.Range("I2").Value = 0.1474788677740264
$ws.Range("J2").Value = 0.1474788677740264
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.110028
$ws.Range("N2").Value = 0.330084
$ws.Range("Q2").Value = 5.195332435052
$ws.Range("R2").Value = 46.757991915468
$ws.Range("S2").Value = 0.1474788677740264
$ws.Range("T2").Value = 0.1474788677740264

# Row 3
$ws.Range("I3").Value = 0.7320002818921112
$ws.Range("J3").Value = 0.7320002818921111
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.110028
$ws.Range("N3").Value = 0.330084
$ws.Range("Q3").Value = 25.786642278868
$ws.Range("R3").Value = 232.079780509812
$ws.Range("S3").Value = 0.7320002818921112
$ws.Range("T3").Value = 0.7320002818921111

# Row 4
$ws.Range("G4").Value = 5.398689
$ws.Range("H4").Value = 16.196067
$ws.Range("I4").Value = 0.01686195715414818
$ws.Range("J4").Value = 0.01686195715414818
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.110028
$ws.Range("N4").Value = 0.330084
$ws.Range("Q4").Value = 0.594006953292
$ws.Range("R4").Value = 5.346062579628
$ws.Range("S4").Value = 0.01686195715414818
$ws.Range("T4").Value = 0.01686195715414818

# Row 5
$ws.Range("G5").Value = 24.22137033333334
$ws.Range("H5").Value = 72.66411100000001
$ws.Range("I5").Value = 0.07565164594134288
$ws.Range("J5").Value = 0.07565164594134288
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.110028
$ws.Range("N5").Value = 0.330084
$ws.Range("Q5").Value = 2.665028935036
$ws.Range("R5").Value = 23.985260415324
$ws.Range("S5").Value = 0.07565164594134288
$ws.Range("T5").Value = 0.07565164594134288

# Row 6
$ws.Range("G6").Value = 1.526601
$ws.Range("H6").Value = 4.579803
$ws.Range("I6").Value = 0.004768098450101454
$ws.Range("J6").Value = 0.004768098450101453
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.110028
$ws.Range("N6").Value = 0.330084
$ws.Range("Q6").Value = 0.167968854828
$ws.Range("R6").Value = 1.511719693452
$ws.Range("S6").Value = 0.004768098450101454
$ws.Range("T6").Value = 0.004768098450101453

# Row 7
$ws.Range("G7").Value = 7.440473
$ws.Range("H7").Value = 22.321419
$ws.Range("I7").Value = 0.02323914878826996
$ws.Range("J7").Value = 0.02323914878826995
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.110028
$ws.Range("N7").Value = 0.330084
$ws.Range("Q7").Value = 0.818660363244
$ws.Range("R7").Value = 7.367943269195999
$ws.Range("S7").Value = 0.02323914878826996
$ws.Range("T7").Value = 0.02323914878826995
